$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H98").Value = 1717.3334
$ws.Range("I98").Value = 1638
$ws.Range("J98").Value = 1955.3334
$ws.Range("K98").Value = 1638
$ws.Range("L98").Value = 1955.3334
$ws.Range("M98").Value = -140
$ws.Range("N98").Value = -4951.3334

$ws.Range("H99").Value = 1008.3333
$ws.Range("I99").Value = 512.5
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1537.5
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -39.5
$ws.Range("N99").Value = -8996

$ws.Range("H101").Value = 20003816
$ws.Range("J101").Value = 728
$ws.Range("L101").Value = 2184
$ws.Range("N101").Value = -5428

$ws.Range("H116").Value = 7249.5
$ws.Range("J116").Value = 7999.6665
$ws.Range("L116").Value = 7999.6665
$ws.Range("N116").Value = -14883.6665

$ws.Range("H122").Value = 1717.3334
$ws.Range("I122").Value = 1638
$ws.Range("J122").Value = 1955.3334
$ws.Range("K122").Value = 4914
$ws.Range("L122").Value = 5866.0002
$ws.Range("M122").Value = -2464
$ws.Range("N122").Value = -10766.0002

$ws.Range("H138").Value = 2248.0908
$ws.Range("I138").Value = 1880.1765
$ws.Range("K138").Value = 5640.529500000001
$ws.Range("M138").Value = -500.5295000000006

$ws.Range("H141").Value = 2900
$ws.Range("I141").Value = 1850
$ws.Range("K141").Value = 5550
$ws.Range("M141").Value = -370

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1983.7368
$ws.Range("I2").Value = 1393.8334
$ws.Range("J2").Value = 2995
$ws.Range("K2").Value = 1393.8334
$ws.Range("L2").Value = 2995
$ws.Range("M2").Value = -1280.8334
$ws.Range("N2").Value = -3221

$ws.Range("H97").Value = 1192
$ws.Range("I97").Value = 1128.8
$ws.Range("K97").Value = 1128.8
$ws.Range("M97").Value = -632.8

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws.Range("H110").Value = 1809.25
$ws.Range("I110").Value = 1946.4546
$ws.Range("K110").Value = 1946.4546
$ws.Range("M110").Value = 98.54539999999997

$ws.Range("H116").Value = 1983.7368
$ws.Range("I116").Value = 1393.8334
$ws.Range("J116").Value = 2995
$ws.Range("K116").Value = 1393.8334
$ws.Range("L116").Value = 2995
$ws.Range("M116").Value = 900.1666
$ws.Range("N116").Value = -7583

$ws.Range("H122").Value = 2362.889
$ws.Range("I122").Value = 1660.2
$ws.Range("J122").Value = 3241.25
$ws.Range("K122").Value = 4980.6
$ws.Range("L122").Value = 9723.75
$ws.Range("M122").Value = -2530.6
$ws.Range("N122").Value = -14623.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1983.7368
$ws.Range("I3").Value = 1393.8334
$ws.Range("J3").Value = 2995
$ws.Range("K3").Value = 1393.8334
$ws.Range("L3").Value = 2995
$ws.Range("M3").Value = -1279.8334
$ws.Range("N3").Value = -3223

$ws.Range("H99").Value = 3999.3333
$ws.Range("J99").Value = 3999.3333
$ws.Range("L99").Value = 3999.3333
$ws.Range("N99").Value = -6995.3333

$ws.Range("H107").Value = 2225.5
$ws.Range("I107").Value = 800.8333
$ws.Range("K107").Value = 800.8333
$ws.Range("M107").Value = 1119.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2393.889
$ws.Range("I31").Value = 2393.889
$ws.Range("K31").Value = 2393.889
$ws.Range("M31").Value = -2098.889

$ws.Range("H34").Value = 2393.889
$ws.Range("I34").Value = 2393.889
$ws.Range("K34").Value = 2393.889
$ws.Range("M34").Value = -2191.889

$ws.Range("H62").Value = 1888.5
$ws.Range("I62").Value = 1888.5
$ws.Range("K62").Value = 1888.5
$ws.Range("M62").Value = -1264.5

$ws.Range("H65").Value = 1888.5
$ws.Range("I65").Value = 1888.5
$ws.Range("K65").Value = 9442.5
$ws.Range("M65").Value = -6322.5

$ws.Range("H134").Value = 2435.625
$ws.Range("I134").Value = 2355
$ws.Range("K134").Value = 7065
$ws.Range("M134").Value = -4530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 107
$ws.Range("I40").Value = 107
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 428
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -359
$ws.Range("N40").ClearContents()

$ws.Range("H76").Value = 9249.833000000001
$ws.Range("I76").Value = 7333
$ws.Range("J76").Value = 11166.667
$ws.Range("K76").Value = 21999
$ws.Range("L76").Value = 33500.001
$ws.Range("M76").Value = -21616
$ws.Range("N76").Value = -34266.001

$ws.Range("H79").Value = 9249.833000000001
$ws.Range("I79").Value = 7333
$ws.Range("J79").Value = 11166.667
$ws.Range("K79").Value = 21999
$ws.Range("L79").Value = 33500.001
$ws.Range("M79").Value = -20673
$ws.Range("N79").Value = -36152.001

$ws.Range("H80").Value = 4444
$ws.Range("J80").Value = 4444
$ws.Range("L80").Value = 13332
$ws.Range("N80").Value = -15204

$ws.Range("H83").Value = 4444
$ws.Range("J83").Value = 4444
$ws.Range("L83").Value = 39996
$ws.Range("N83").Value = -49356

$ws.Range("H92").Value = 525.25
$ws.Range("I92").Value = 525.25
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1575.75
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -327.75
$ws.Range("N92").ClearContents()

$ws.Range("H94").Value = 8000
$ws.Range("J94").Value = 15000
$ws.Range("L94").Value = 45000
$ws.Range("N94").Value = -46352

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws.Range("H115").Value = 4000
$ws.Range("J115").Value = 4000
$ws.Range("L115").Value = 12000
$ws.Range("N115").Value = -14350

$ws.Range("H118").Value = 2809
$ws.Range("I118").Value = 2213.5
$ws.Range("K118").Value = 6640.5
$ws.Range("M118").Value = -5397.5

$ws.Range("H119").Value = 900
$ws.Range("I119").Value = 900
$ws.Range("K119").Value = 2700
$ws.Range("M119").Value = 2138

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16684.428
$ws.Range("I80").Value = 2859
$ws.Range("J80").Value = 51248
$ws.Range("K80").Value = 2859
$ws.Range("L80").Value = 51248
$ws.Range("M80").Value = -1861
$ws.Range("N80").Value = -53244

$ws.Range("H83").Value = 16684.428
$ws.Range("I83").Value = 2859
$ws.Range("J83").Value = 51248
$ws.Range("K83").Value = 14295
$ws.Range("L83").Value = 256240
$ws.Range("M83").Value = -9303
$ws.Range("N83").Value = -266224

$ws.Range("H102").Value = 2128.4443
$ws.Range("I102").Value = 2144.5
$ws.Range("K102").Value = 2144.5
$ws.Range("M102").Value = -522.5

$ws.Range("H107").Value = 1616.6666
$ws.Range("J107").Value = 2932.25
$ws.Range("L107").Value = 2932.25
$ws.Range("N107").Value = -6772.25

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6714.8184
$ws.Range("J7").Value = 7551.4443
$ws.Range("L7").Value = 7551.4443
$ws.Range("N7").Value = -7775.4443

$ws.Range("H22").Value = 1635.1538
$ws.Range("I22").Value = 1872.2858
$ws.Range("J22").Value = 1358.5
$ws.Range("K22").Value = 1872.2858
$ws.Range("L22").Value = 1358.5
$ws.Range("M22").Value = -1577.2858
$ws.Range("N22").Value = -1948.5

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()

$ws.Range("H27").Value = 1635.1538
$ws.Range("I27").Value = 1872.2858
$ws.Range("J27").Value = 1358.5
$ws.Range("K27").Value = 1872.2858
$ws.Range("L27").Value = 1358.5
$ws.Range("M27").Value = -1765.2858
$ws.Range("N27").Value = -1572.5

$ws.Range("H40").Value = 3705.25
$ws.Range("I40").Value = 3408
$ws.Range("K40").Value = 3408
$ws.Range("M40").Value = -3272

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H126").Value = 6714.8184
$ws.Range("J126").Value = 7551.4443
$ws.Range("L126").Value = 22654.3329
$ws.Range("N126").Value = -27594.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 28444
$ws.Range("I62").Value = 28000
$ws.Range("J62").Value = 28888
$ws.Range("K62").Value = 28000
$ws.Range("L62").Value = 28888
$ws.Range("M62").Value = -27376
$ws.Range("N62").Value = -30136

$ws.Range("H65").Value = 28444
$ws.Range("I65").Value = 28000
$ws.Range("J65").Value = 28888
$ws.Range("K65").Value = 140000
$ws.Range("L65").Value = 144440
$ws.Range("M65").Value = -136880
$ws.Range("N65").Value = -150680

$ws.Range("H107").Value = 2038.9375
$ws.Range("I107").Value = 1826.8334
$ws.Range("J107").Value = 2675.25
$ws.Range("K107").Value = 5480.5002
$ws.Range("L107").Value = 8025.75
$ws.Range("M107").Value = -3560.5002
$ws.Range("N107").Value = -11865.75

$ws.Range("H113").Value = 506
$ws.Range("I113").Value = 387.5
$ws.Range("J113").Value = 743
$ws.Range("K113").Value = 1162.5
$ws.Range("L113").Value = 2229
$ws.Range("M113").Value = 1007.5
$ws.Range("N113").Value = -6569
